# Update countries & provincias Spain
#
# Source data refresh: the "Datos actualizados" timestamp moves from 14:20
# to 14:50, several provinces get refreshed case counts, and three pairs of
# rows swap places in sort order (Araba/Alava now outranks Alacant/Alicante,
# Gipuzkoa/Guipuzcoa jumps above Cantabria/Granada/Caceres, and Lanzarote
# overtakes Melilla). Rather than physically re-sorting rows, each affected
# cell is written explicitly with its final value so every row keeps the
# formatting/style already on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 14:50"

# row -> @{ col = value }  (only columns that actually change are listed)
$updates = [ordered]@{
    6  = @{ B = 3682; C = 2490; D = 4415; E = 187 }                      # Bizkaia/Vizcaya
    11 = @{ A = "Araba/Alava";        B = 2347; C = 2490; D = 4415; E = 167 }
    12 = @{ A = "Alacant/Alicante";   B = 2320; C = 172;  D = 1960; E = 188 }
    20 = @{ C = 94; D = 1282 }                                           # Tenerife
    24 = @{ A = "Gipuzkoa/Guipuzcoa"; B = 1288; C = 2490; D = 4415; E = 58 }
    25 = @{ A = "Cantabria";          B = 1268; C = 60;   D = 1148; E = 60 }
    26 = @{ A = "Granada";            B = 1230; C = 15;   D = 1129; E = 86 }
    27 = @{ A = "Caceres";            B = 1212; C = 45;   D = 1012; E = 155 }
    44 = @{ B = 396; C = 94; D = 1282; E = 19 }                          # Gran Canaria
    53 = @{ B = 69;  C = 94; D = 1282; E = 2 }                           # La Palma
    54 = @{ A = "Lanzarote"; B = 65; C = 94; D = 1282; E = 3 }
    55 = @{ A = "Melilla";   B = 62; C = 0;  D = 61;   E = 1 }
    58 = @{ B = 33; C = 94; D = 1282; E = 0 }                            # Fuerteventura
    61 = @{ C = 94; D = 1282; E = 0 }                                    # La Gomera
    63 = @{ C = 94; D = 1282; E = 0 }                                    # El Hierro
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
